$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Quarterly update: insert two new columns (two new reporting quarters) right
# before column D. This shifts the existing quarter columns D:K -> F:M, and
# leaves two empty columns D:E ready to receive the newest quarter figures.
# ---------------------------------------------------------------------------
$ws.Columns("D:E").Insert()

# Copy the number formatting / style that the (now shifted) first data column
# F carries into the two freshly inserted columns D:E, so the new cells look
# identical to the rest of the quarterly columns (date style for the header
# row, numeric style for the data rows). Work block-by-block so that the
# section-heading rows that have no data cells (36/37, 78/79) are not
# touched and stay free of any stray D/E cells.
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)

$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)

$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Populate the two new quarter columns with the latest reported figures.
# ---------------------------------------------------------------------------

# Income Statement - Period Ending
$ws.Range("D7").Value = 43496
$ws.Range("E7").Value = 43404

# Total Revenue
$ws.Range("D8").Value = 6400
$ws.Range("E8").Value = 6700

# Cost of Revenue
$ws.Range("D9").Value = 3400
$ws.Range("E9").Value = 3500

# Gross Profit
$ws.Range("D10").Value = 3000
$ws.Range("E10").Value = 3200

# Research Development
$ws.Range("D12").Value = 1300
$ws.Range("E12").Value = 1200

# Selling General and Administrative
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0

# Non Recurring
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0

# Others
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0

# Total Operating Expenses
$ws.Range("D17").Value = 6800
$ws.Range("E17").Value = 6400

# Operating Income or Loss
$ws.Range("D18").Value = -400
$ws.Range("E18").Value = 300

# Total Other Income/Expenses Net
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0

# Earnings Before Interest And Taxes
$ws.Range("D21").Value = -200
$ws.Range("E21").Value = 400

# Interest Expense
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0

# Income Before Tax
$ws.Range("D23").Value = -400
$ws.Range("E23").Value = 300

# Income Tax Expense
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0

# Minority Interest
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0

# Income After Tax
$ws.Range("D26").Value = -400
$ws.Range("E26").Value = 300

# Net Income From Continuing Ops
$ws.Range("D27").Value = -400
$ws.Range("E27").Value = 300

# Non-recurring Events
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0

# Discontinued Operations
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0

# Extraordinary Items
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0

# Effect Of Accounting Changes
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0

# Other Items
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0

# Net Income
$ws.Range("D33").Value = -400
$ws.Range("E33").Value = 300

# Preferred Stock And Other Adjustments
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0

# Net Income Applicable To Common Shares
$ws.Range("D35").Value = -400
$ws.Range("E35").Value = 300

# Balance Sheet - Period Ending
$ws.Range("D38").Value = 43496
$ws.Range("E38").Value = 43404

# Cash And Cash Equivalents
$ws.Range("D41").Value = 3300
$ws.Range("E41").Value = 2000

# Short Term Investments
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0

# Net Receivables
$ws.Range("D43").Value = 3900
$ws.Range("E43").Value = 4700

# Inventory
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0

# Other Current Assets
$ws.Range("D45").Value = 300
$ws.Range("E45").Value = 300

# Total Current Assets
$ws.Range("D46").Value = 7500
$ws.Range("E46").Value = 6900

# Long Term Investments
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0

# Property Plant and Equipment
$ws.Range("D48").Value = 2600
$ws.Range("E48").Value = 2500

# Goodwill
$ws.Range("D49").Value = 700
$ws.Range("E49").Value = 700

# Intangible Assets
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0

# Accumulated Amortization
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0

# Other Assets
$ws.Range("D52").Value = 100
$ws.Range("E52").Value = 100

# Deferred Long Term Asset Charges
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0

# Total Assets
$ws.Range("D54").Value = 10900
$ws.Range("E54").Value = 10200

# Accounts Payable
$ws.Range("D57").Value = 2100
$ws.Range("E57").Value = 2000

# Short/Current Long Term Debt
$ws.Range("D58").Value = 200
$ws.Range("E58").Value = 100

# Other Current Liabilities
$ws.Range("D59").Value = 5400
$ws.Range("E59").Value = 5200

# Total Current Liabilities
$ws.Range("D60").Value = 7700
$ws.Range("E60").Value = 7300

# Long Term Debt
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 100

# Other Liabilities
$ws.Range("D62").Value = 900
$ws.Range("E62").Value = 800

# Deferred Long Term Liability Charges
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0

# Negative Goodwill
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0

# Total Liabilities
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0

# Misc Stocks Options Warrants
$ws.Range("D66").Value = 8600
$ws.Range("E66").Value = 8200

# Redeemable Preferred Stock
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0

# Preferred Stock
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0

# Common Stock
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0

# Retained Earnings
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0

# Treasury Stock
$ws.Range("D72").Value = -70400
$ws.Range("E72").Value = -70100

# Capital Surplus
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0

# Other Stockholder Equity
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0

# Total Stockholder Equity
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0

# Net Tangible Assets
$ws.Range("D76").Value = 2300
$ws.Range("E76").Value = 2000

# (unlabeled trailing Balance Sheet row)
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0

# Cash Flow Statement - Period Ending
$ws.Range("D80").Value = 43496
$ws.Range("E80").Value = 43404

# Net Income
$ws.Range("D81").Value = -400
$ws.Range("E81").Value = 300

# Depreciation
$ws.Range("D83").Value = 200
$ws.Range("E83").Value = 200

# Adjustments To Net Income
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0

# Changes In Accounts Receivables
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0

# Changes In Liabilities
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0

# Changes In Inventories
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0

# Changes In Other Operating Activities
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0

# Total Cash Flow From Operating Activities
$ws.Range("D89").Value = 1200
$ws.Range("E89").Value = 100

# Capital Expenditures
$ws.Range("D91").Value = -200
$ws.Range("E91").Value = -300

# Investments
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0

# Other Cashflows from Investing Activities
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0

# Total Cash Flows From Investing Activities
$ws.Range("D94").Value = -200
$ws.Range("E94").Value = -300

# Dividends Paid
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0

# Sale Purchase of Stock
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0

# Net Borrowings
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0

# Other Cash Flows from Financing Activities
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0

# Total Cash Flows From Financing Activities
$ws.Range("D100").Value = 400
$ws.Range("E100").Value = 1000

# Effect Of Exchange Rate Changes
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0

# Change In Cash and Cash Equivalents
$ws.Range("D102").Value = 1400
$ws.Range("E102").Value = 800
